$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.098.25'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.830.56'
$ws.Range('E3').Value = '  -0.33%  '
$c = $ws.Range('D4')
$c.Value = "'0.9994"
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$c = $ws.Range('D5')
$c.Value = "'243.15"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.11%  '
$c = $ws.Range('D6')
$c.Value = "'0.6256"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E7').Value = '  -0.14%  '
$c = $ws.Range('D8')
$c.Value = "'0.07499"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -1.18%  '
$c = $ws.Range('D9')
$c.Value = "'0.2925"
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.09%  '
$c = $ws.Range('D10')
$c.Value = "'23.17"
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +2.46%  '
$c = $ws.Range('D11')
$c.Value = "'0.07681"
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('D12').Value = '1.827.25'
$ws.Range('E12').Value = '  -0.86%  '
$c = $ws.Range('D13')
$c.Value = "'5.021"
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +0.95%  '
$c = $ws.Range('D14')
$c.Value = "'0.6677"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.33%  '
$c = $ws.Range('D15')
$c.Value = "'82.75"
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.42%  '
$c = $ws.Range('D16')
$c.Value = "'0.000009383"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -6.28%  '
$c = $ws.Range('D17')
$c.Value = "'5.986"
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('D18').Value = '29.094.59'
$ws.Range('D19').Value = '2.077.63'
$ws.Range('E19').Value = '  -0.60%  '
$c = $ws.Range('D20')
$c.Value = "'12.60"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +1.44%  '
$c = $ws.Range('D21')
$c.Value = "'222.78"
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.89%  '
$ws.Range('E22').Value = '  +0.06%  '
$c = $ws.Range('D23')
$c.Value = "'7.144"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -1.10%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D24')
$c.Value = "'1.001"
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D25')
$c.Value = "'160.25"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D26')
$c.Value = "'0.1394"
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D27')
$c.Value = "'8.490"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D28')
$c.Value = "'17.90"
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D29')
$c.Value = "'1.494"
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D30')
$c.Value = "'0.05824"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +10.74%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D31')
$c.Value = "'4.154"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D32')
$c.Value = "'4.116"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +2.27%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D33')
$c.Value = "'1.208"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D34')
$c.Value = "'0.7411"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.76%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D35')
$c.Value = "'1.828"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.87%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D36')
$c.Value = "'1.139"
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D37')
$c.Value = "'2.668"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.30%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.228.03'
$ws.Range('E38').Value = '  -1.06%  '
$c = $ws.Range('D39')
$c.Value = "'2.763"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D40')
$c.Value = "'0.01778"
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.38%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D41')
$c.Value = "'6.486"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D42')
$c.Value = "'0.8892"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range('D43')
$c.Value = "'1.001"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D44')
$c.Value = "'102.07"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.978.57'
$ws.Range('E45').Value = '  -0.39%  '
$c = $ws.Range('D46')
$c.Value = "'0.00000000125"
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D47')
$c.Value = "'65.91"
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +2.31%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D48')
$c.Value = "'0.5087"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('B49').Value = 'XinFinNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$c = $ws.Range('D49')
$c.Value = "'0.07559"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +12.79%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range('D50')
$c.Value = "'0.4062"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D51')
$c.Value = "'8.955"
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +1.00%  '
